$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/value cells (safe to assign directly; Excel will not
#     mis-parse these as numbers because of extra dots, %, spaces, etc.) ---
$ws.Range("D2").Value = '61.462.08'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '2.381.62'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '2.382.95'
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("E10").Value = '  +3.13%  '
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("E12").Value = '  +1.62%  '
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("E15").Value = '  +2.90%  '
$ws.Range("D16").Value = '61.367.71'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").Value = '2.377.21'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("E18").Value = '  +3.18%  '
$ws.Range("E19").Value = '  +1.85%  '
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("E24").Value = '  -8.79%  '
$ws.Range("E25").Value = '  +6.92%  '
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = '0.0₃0899'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("E29").Value = '  +3.74%  '
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  +4.56%  '
$ws.Range("E35").Value = '  +1.86%  '
$ws.Range("E36").Value = '  +2.37%  '
$ws.Range("E37").Value = '  +1.86%  '
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("E39").Value = '  +4.76%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  +2.33%  '
$ws.Range("E42").Value = '  +7.98%  '
$ws.Range("E43").Value = '  +2.06%  '
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("E45").Value = '  +1.82%  '
$ws.Range("E46").Value = '  -1.23%  '
$ws.Range("E47").Value = '  +1.63%  '
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("E51").Value = '  +0.93%  '

# --- Numeric-looking price strings that must stay literal text (they
#     carry significant trailing zeros, e.g. "25.30"). Assign each as a
#     formula producing the literal string, then convert that single
#     cell's formula to a static value via copy / paste-special values so
#     Excel does not re-parse the text as a number and strip the zero. ---
$ws.Range("D5").Formula = '="550.52"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="139.67"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D8").Formula = '="0.524"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("D14").Formula = '="25.30"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D18").Formula = '="10.98"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("D19").Formula = '="321.59"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="6.78"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D23").Formula = '="64.47"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D24").Formula = '="1.70"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D25").Formula = '="8.64"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D26").Formula = '="8.19"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D27").Formula = '="515.51"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D30").Formula = '="1.38"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D36").Formula = '="5.49"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="18.56"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D39").Formula = '="146.61"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="41.21"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="150.47"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D43").Formula = '="2.15"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("D46").Formula = '="19.50"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("D51").Formula = '="16.83"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$excel.CutCopyMode = 0

Write-Host "Updated cryptos list"
